$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("latest")

$ws.Cells.Item(2, 2).Value = 0.5553305257193456
$ws.Cells.Item(2, 3).Value = 0.6136273584998631
$ws.Cells.Item(2, 4).Value = 0.5228536874213209
$ws.Cells.Item(2, 5).Value = 0.723086224057215
$ws.Cells.Item(2, 6).Value = 0.4805816767817249
$ws.Cells.Item(2, 7).Value = 14

$ws.Cells.Item(3, 2).Value = 0.3621708991783025
$ws.Cells.Item(3, 3).Value = 0.4983652827815714
$ws.Cells.Item(3, 4).Value = 0.342784039812318
$ws.Cells.Item(3, 5).Value = 0.585477616832888
$ws.Cells.Item(3, 6).Value = 0.4788015972899659
$ws.Cells.Item(3, 7).Value = 13

$ws.Cells.Item(4, 2).Value = 0.322109523816475
$ws.Cells.Item(4, 3).Value = 0.4916759303336116
$ws.Cells.Item(4, 4).Value = 0.3220049964818212
$ws.Cells.Item(4, 5).Value = 0.5674548409184833
$ws.Cells.Item(4, 6).Value = 0.4879461048650333
$ws.Cells.Item(4, 7).Value = 12

$ws.Cells.Item(5, 2).Value = 0.4479694299613027
$ws.Cells.Item(5, 3).Value = 0.5810026869938167
$ws.Cells.Item(5, 4).Value = 0.4406117130697931
$ws.Cells.Item(5, 5).Value = 0.6637858939972987
$ws.Cells.Item(5, 6).Value = 0.5137398302438039
$ws.Cells.Item(5, 7).Value = 11

$ws.Cells.Item(6, 2).Value = 0.4056295499570255
$ws.Cells.Item(6, 3).Value = 0.5539856660635349
$ws.Cells.Item(6, 4).Value = 0.4125774920458339
$ws.Cells.Item(6, 5).Value = 0.6423219535761128
$ws.Cells.Item(6, 6).Value = 0.5249784760111545
$ws.Cells.Item(6, 7).Value = 10

$ws.Cells.Item(7, 2).Value = 0.3136617854706863
$ws.Cells.Item(7, 3).Value = 0.4838333121955758
$ws.Cells.Item(7, 4).Value = 0.3116578194982599
$ws.Cells.Item(7, 5).Value = 0.5582632170385757
$ws.Cells.Item(7, 6).Value = 0.4898299366237236
$ws.Cells.Item(7, 7).Value = 9

$ws.Cells.Item(8, 2).Value = 0.3613000660075406
$ws.Cells.Item(8, 3).Value = 0.5407655076793962
$ws.Cells.Item(8, 4).Value = 0.3759763669136598
$ws.Cells.Item(8, 5).Value = 0.6131691177103261
$ws.Cells.Item(8, 6).Value = 0.5296237254251971
$ws.Cells.Item(8, 7).Value = 8

$ws.Cells.Item(9, 2).Value = 0.5028015250901602
$ws.Cells.Item(9, 3).Value = 0.5533246075941491
$ws.Cells.Item(9, 4).Value = 0.4051903875053499
$ws.Cells.Item(9, 5).Value = 0.6365456680438175
$ws.Cells.Item(9, 6).Value = 0.4216371064289196
$ws.Cells.Item(9, 7).Value = 7

$ws.Cells.Item(10, 2).Value = 0.4053476105038098
$ws.Cells.Item(10, 3).Value = 0.459336539820115
$ws.Cells.Item(10, 4).Value = 0.2581058573136952
$ws.Cells.Item(10, 5).Value = 0.5080411964729782
$ws.Cells.Item(10, 6).Value = 0.33549814659258
$ws.Cells.Item(10, 7).Value = 6

$ws.Cells.Item(11, 2).Value = 0.4015449179333902
$ws.Cells.Item(11, 3).Value = 0.4597420797009394
$ws.Cells.Item(11, 4).Value = 0.2670924375257908
$ws.Cells.Item(11, 5).Value = 0.516809865933102
$ws.Cells.Item(11, 6).Value = 0.3637549250657264
$ws.Cells.Item(11, 7).Value = 5

$wb.Save()
